$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 427.08334
$ws.Range("I2").Value = 286.14285
$ws.Range("K2").Value = 286.14285
$ws.Range("M2").Value = -173.14285
$ws.Range("H3").Value = 57250
$ws.Range("J3").Value = 57250
$ws.Range("L3").Value = 57250
$ws.Range("N3").Value = -57478
$ws.Range("H17").Value = 2017.8334
$ws.Range("J17").Value = 2017.8334
$ws.Range("L17").Value = 6053.5002
$ws.Range("N17").Value = -6389.5002
$ws.Range("H40").Value = 3954.524
$ws.Range("I40").Value = 4798.25
$ws.Range("J40").Value = 3435.3076
$ws.Range("K40").Value = 4798.25
$ws.Range("L40").Value = 3435.3076
$ws.Range("M40").Value = -4623.25
$ws.Range("N40").Value = -3785.3076
$ws.Range("H51").Value = 8863.666999999999
$ws.Range("I51").Value = 9300
$ws.Range("J51").Value = 8645.5
$ws.Range("K51").Value = 9300
$ws.Range("L51").Value = 8645.5
$ws.Range("M51").Value = -8816
$ws.Range("N51").Value = -9613.5
$ws.Range("H55").Value = 517.1667
$ws.Range("J55").Value = 199
$ws.Range("L55").Value = 199
$ws.Range("N55").Value = -627
$ws.Range("H94").Value = 40000
$ws.Range("I94").Value = 40000
$ws.Range("K94").Value = 40000
$ws.Range("M94").Value = -39549
$ws.Range("H98").Value = 1863.9149
$ws.Range("I98").Value = 1866.7778
$ws.Range("K98").Value = 1866.7778
$ws.Range("M98").Value = -368.7778000000001
$ws.Range("H102").Value = 57250
$ws.Range("J102").Value = 57250
$ws.Range("L102").Value = 57250
$ws.Range("N102").Value = -63740
$ws.Range("H103").Value = 3828.6
$ws.Range("I103").Value = 4410.75
$ws.Range("J103").Value = 1500
$ws.Range("K103").Value = 13232.25
$ws.Range("L103").Value = 4500
$ws.Range("M103").Value = -12646.25
$ws.Range("N103").Value = -5672
$ws.Range("H122").Value = 1863.9149
$ws.Range("I122").Value = 1866.7778
$ws.Range("K122").Value = 5600.3334
$ws.Range("M122").Value = -3150.3334
$ws.Range("H135").Value = 1289.8334
$ws.Range("I135").Value = 940.8333
$ws.Range("K135").Value = 8467.4997
$ws.Range("M135").Value = -5932.4997
$ws.Range("H138").Value = 365672.56
$ws.Range("J138").Value = 598509.25
$ws.Range("L138").Value = 1795527.75
$ws.Range("N138").Value = -1805807.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4062.8572
$ws.Range("I63").Value = 4062.8572
$ws.Range("K63").Value = 4062.8572
$ws.Range("M63").Value = -3376.8572
$ws.Range("H66").Value = 4062.8572
$ws.Range("I66").Value = 4062.8572
$ws.Range("K66").Value = 20314.286
$ws.Range("M66").Value = -16882.286
$ws.Range("H102").Value = 3027.4358
$ws.Range("I102").Value = 2414.8125
$ws.Range("K102").Value = 2414.8125
$ws.Range("M102").Value = -792.8125
$ws.Range("H112").Value = 29999.5
$ws.Range("J112").Value = 29999.5
$ws.Range("L112").Value = 29999.5
$ws.Range("N112").Value = -32953.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 11820872
$ws.Range("I105").Value = 716192.4399999999
$ws.Range("J105").Value = 31254062
$ws.Range("K105").Value = 716192.4399999999
$ws.Range("L105").Value = 31254062
$ws.Range("M105").Value = -714445.4399999999
$ws.Range("N105").Value = -31257556
$ws.Range("H107").Value = 1501.0476
$ws.Range("I107").Value = 1436.875
$ws.Range("J107").Value = 1706.4
$ws.Range("K107").Value = 1436.875
$ws.Range("L107").Value = 1706.4
$ws.Range("M107").Value = 483.125
$ws.Range("N107").Value = -5546.4
$ws.Range("H134").Value = 2324.5881
$ws.Range("I134").Value = 2051.48
$ws.Range("K134").Value = 6154.440000000001
$ws.Range("M134").Value = -3619.440000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4664.8057
$ws.Range("I31").Value = 3533.5417
$ws.Range("J31").Value = 6927.3335
$ws.Range("K31").Value = 3533.5417
$ws.Range("L31").Value = 6927.3335
$ws.Range("M31").Value = -3238.5417
$ws.Range("N31").Value = -7517.3335
$ws.Range("H34").Value = 4664.8057
$ws.Range("I34").Value = 3533.5417
$ws.Range("J34").Value = 6927.3335
$ws.Range("K34").Value = 3533.5417
$ws.Range("L34").Value = 6927.3335
$ws.Range("M34").Value = -3331.5417
$ws.Range("N34").Value = -7331.3335
$ws.Range("H42").Value = 27833
$ws.Range("J42").Value = 27833
$ws.Range("L42").Value = 27833
$ws.Range("N42").Value = -29019
$ws.Range("H62").Value = 11116605
$ws.Range("I62").Value = 20005200
$ws.Range("J62").Value = 5862.25
$ws.Range("K62").Value = 20005200
$ws.Range("L62").Value = 5862.25
$ws.Range("M62").Value = -20004576
$ws.Range("N62").Value = -7110.25
$ws.Range("H65").Value = 11116605
$ws.Range("I65").Value = 20005200
$ws.Range("J65").Value = 5862.25
$ws.Range("K65").Value = 100026000
$ws.Range("L65").Value = 29311.25
$ws.Range("M65").Value = -100022880
$ws.Range("N65").Value = -35551.25
$ws.Range("H134").Value = 4539.4736
$ws.Range("I134").Value = 4328.125
$ws.Range("J134").Value = 5666.6665
$ws.Range("K134").Value = 12984.375
$ws.Range("L134").Value = 16999.9995
$ws.Range("M134").Value = -10449.375
$ws.Range("N134").Value = -22069.9995

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 2114
$ws.Range("I2").Value = 37
$ws.Range("J2").Value = 3844.8333
$ws.Range("K2").Value = 222
$ws.Range("L2").Value = 23068.9998
$ws.Range("M2").Value = -109
$ws.Range("N2").Value = -23294.9998
$ws.Range("H97").Value = 833787
$ws.Range("J97").Value = 560.5
$ws.Range("L97").Value = 1681.5
$ws.Range("N97").Value = -2673.5
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("H131").Value = 2746.1
$ws.Range("I131").Value = 4032.6667
$ws.Range("J131").Value = 2194.7144
$ws.Range("K131").Value = 12098.0001
$ws.Range("L131").Value = 6584.1432
$ws.Range("M131").Value = -7058.000100000001
$ws.Range("N131").Value = -16664.1432
$ws.Range("H138").Value = 4847.8335
$ws.Range("I138").Value = 3055.25
$ws.Range("K138").Value = 9165.75
$ws.Range("M138").Value = -4025.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 711.5
$ws.Range("J2").Value = 205.5
$ws.Range("L2").Value = 205.5
$ws.Range("N2").Value = -431.5
$ws.Range("H70").Value = 129908.94
$ws.Range("I70").Value = 171215.42
$ws.Range("J70").Value = 5989.5
$ws.Range("K70").Value = 171215.42
$ws.Range("L70").Value = 5989.5
$ws.Range("M70").Value = -170945.42
$ws.Range("N70").Value = -6529.5
$ws.Range("H73").Value = 129908.94
$ws.Range("I73").Value = 171215.42
$ws.Range("J73").Value = 5989.5
$ws.Range("K73").Value = 171215.42
$ws.Range("L73").Value = 5989.5
$ws.Range("M73").Value = -170279.42
$ws.Range("N73").Value = -7861.5
$ws.Range("H107").Value = 828.3333
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 828.3333
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 828.3333
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -4668.3333
$ws.Range("H113").Value = 3602.182
$ws.Range("I113").Value = 2961.1
$ws.Range("J113").Value = 10013
$ws.Range("K113").Value = 2961.1
$ws.Range("L113").Value = 10013
$ws.Range("M113").Value = -791.0999999999999
$ws.Range("N113").Value = -14353
$ws.Range("H122").Value = 4761.98
$ws.Range("I122").Value = 4632.5
$ws.Range("J122").Value = 5037.125
$ws.Range("K122").Value = 13897.5
$ws.Range("L122").Value = 15111.375
$ws.Range("M122").Value = -11447.5
$ws.Range("N122").Value = -20011.375

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4626.4185
$ws.Range("I40").Value = 4710.029
$ws.Range("J40").Value = 4260.625
$ws.Range("K40").Value = 4710.029
$ws.Range("L40").Value = 4260.625
$ws.Range("M40").Value = -4574.029
$ws.Range("N40").Value = -4532.625
$ws.Range("H132").Value = 7417.9546
$ws.Range("I132").Value = 3707.3845
$ws.Range("K132").Value = 11122.1535
$ws.Range("M132").Value = -8592.1535
$ws.Range("H133").Value = 101325
$ws.Range("J133").Value = 101325
$ws.Range("L133").Value = 101325
$ws.Range("N133").Value = -106385
$ws.Range("H136").Value = 3137.3635
$ws.Range("I136").Value = 2829.818
$ws.Range("K136").Value = 8489.454000000002
$ws.Range("M136").Value = -5939.454000000002

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2937.125
$ws.Range("I81").Value = 2456.7144
$ws.Range("K81").Value = 4913.4288
$ws.Range("M81").Value = -3852.4288
$ws.Range("H84").Value = 2937.125
$ws.Range("I84").Value = 2456.7144
$ws.Range("K84").Value = 24567.144
$ws.Range("M84").Value = -19263.144
$ws.Range("H96").Value = 4729.933
$ws.Range("I96").Value = 4658.1
$ws.Range("J96").Value = 4873.6
$ws.Range("K96").Value = 4658.1
$ws.Range("L96").Value = 4873.6
$ws.Range("M96").Value = -3285.1
$ws.Range("N96").Value = -7619.6
$ws.Range("H113").Value = 458.04166
$ws.Range("J113").Value = 368.55554
$ws.Range("L113").Value = 1105.66662
$ws.Range("N113").Value = -5445.66662
$ws.Range("H122").Value = 11367655
$ws.Range("I122").Value = 4475.294
$ws.Range("J122").Value = 50002464
$ws.Range("K122").Value = 13425.882
$ws.Range("L122").Value = 150007392
$ws.Range("M122").Value = -10975.882
$ws.Range("N122").Value = -150012292
$ws.Range("H126").Value = 2164.6667
$ws.Range("I126").Value = 2047.25
$ws.Range("J126").Value = 2399.5
$ws.Range("K126").Value = 6141.75
$ws.Range("L126").Value = 7198.5
$ws.Range("M126").Value = -3671.75
$ws.Range("N126").Value = -12138.5
$ws.Range("H132").Value = 3329.1924
$ws.Range("I132").Value = 3263.4348
$ws.Range("K132").Value = 9790.304400000001
$ws.Range("M132").Value = -7260.304400000001
